# Add real groups functionality.
# Appends a new "Volleyball" row (row 6) to the Schedule sheet, matching
# the other schedule rows (Group, Date, Time columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

$ws.Cells.Item(6, 1).Value = "Volleyball"
$ws.Cells.Item(6, 2).Value = "okokokok"

# "2021-11-23" looks like a date, so force the cell to text first (via a
# temporary Text number format) to stop it being auto-converted into a
# date serial number, then drop the format change again so the cell is
# left with the sheet's ordinary default formatting (same as the other
# rows).
$ws.Cells.Item(6, 3).NumberFormat = "@"
$ws.Cells.Item(6, 3).Value = "2021-11-23"
$ws.Cells.Item(6, 3).ClearFormats()

$ws.Cells.Item(6, 4).Value = "11:00:00"
